$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 76923576
$ws.Range("I33").Value = 76923576
$ws.Range("K33").Value = 76923576
$ws.Range("M33").Value = -76923347
$ws.Range("H53").Value = 2078.818
$ws.Range("I53").Value = 3743.5
$ws.Range("J53").Value = 81.2
$ws.Range("K53").Value = 3743.5
$ws.Range("L53").Value = 81.2
$ws.Range("M53").Value = -3106.5
$ws.Range("N53").Value = -1355.2
$ws.Range("H131").Value = 4418.125
$ws.Range("I131").Value = 517.25
$ws.Range("K131").Value = 1551.75
$ws.Range("M131").Value = 3488.25
$ws.Range("H137").Value = 1310.1852
$ws.Range("I137").Value = 950.7895
$ws.Range("J137").Value = 2163.75
$ws.Range("K137").Value = 2852.3685
$ws.Range("L137").Value = 6491.25
$ws.Range("M137").Value = -302.3685
$ws.Range("N137").Value = -11591.25
$ws.Range("H138").Value = 2490.3447
$ws.Range("I138").Value = 2047.8572
$ws.Range("J138").Value = 2903.3333
$ws.Range("K138").Value = 6143.571599999999
$ws.Range("L138").Value = 8709.999899999999
$ws.Range("M138").Value = -1003.571599999999
$ws.Range("N138").Value = -18989.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2605.4517
$ws.Range("I61").Value = 2629
$ws.Range("J61").Value = 2264
$ws.Range("K61").Value = 2629
$ws.Range("L61").Value = 2264
$ws.Range("M61").Value = -2417
$ws.Range("N61").Value = -2688
$ws.Range("H74").Value = 3669.1482
$ws.Range("J74").Value = 7048.5557
$ws.Range("L74").Value = 7048.5557
$ws.Range("N74").Value = -8796.555700000001
$ws.Range("H77").Value = 3669.1482
$ws.Range("J77").Value = 7048.5557
$ws.Range("L77").Value = 35242.7785
$ws.Range("N77").Value = -43978.7785
$ws.Range("H132").Value = 3732817
$ws.Range("I132").Value = 4311551.5
$ws.Range("K132").Value = 12934654.5
$ws.Range("M132").Value = -12932124.5
$ws.Range("H136").Value = 2605.4517
$ws.Range("I136").Value = 2629
$ws.Range("J136").Value = 2264
$ws.Range("K136").Value = 7887
$ws.Range("L136").Value = 6792
$ws.Range("M136").Value = -5337
$ws.Range("N136").Value = -11892

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2749.5881
$ws.Range("I134").Value = 1559.871
$ws.Range("J134").Value = 4593.65
$ws.Range("K134").Value = 4679.613
$ws.Range("L134").Value = 13780.95
$ws.Range("M134").Value = -2144.613
$ws.Range("N134").Value = -18850.95

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11167.077
$ws.Range("I31").Value = 13917.2
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 13917.2
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -13622.2
$ws.Range("N31").Value = -2590
$ws.Range("H34").Value = 11167.077
$ws.Range("I34").Value = 13917.2
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 13917.2
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -13715.2
$ws.Range("N34").Value = -2404
$ws.Range("H58").Value = 3054.077
$ws.Range("I58").Value = 1062.25
$ws.Range("J58").Value = 3939.3333
$ws.Range("K58").Value = 1062.25
$ws.Range("L58").Value = 3939.3333
$ws.Range("M58").Value = -859.25
$ws.Range("N58").Value = -4345.3333
$ws.Range("H62").Value = 2990
$ws.Range("H65").Value = 2990
$ws.Range("H132").Value = 12849.182
$ws.Range("I132").Value = 18902.666
$ws.Range("J132").Value = 5585
$ws.Range("K132").Value = 56707.99800000001
$ws.Range("L132").Value = 16755
$ws.Range("M132").Value = -54177.99800000001
$ws.Range("N132").Value = -21815
$ws.Range("H134").Value = 2419.8235
$ws.Range("I134").Value = 2258.5625
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 6775.6875
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -4240.6875
$ws.Range("N134").Value = -20070
$ws.Range("H136").Value = 3054.077
$ws.Range("I136").Value = 1062.25
$ws.Range("J136").Value = 3939.3333
$ws.Range("K136").Value = 3186.75
$ws.Range("L136").Value = 11817.9999
$ws.Range("M136").Value = -636.75
$ws.Range("N136").Value = -16917.9999
$ws.Range("H140").Value = 27271.428
$ws.Range("J140").Value = 27271.428
$ws.Range("L140").Value = 27271.428
$ws.Range("N140").Value = -37631.428

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3906.2058
$ws.Range("I126").Value = 2487.6
$ws.Range("J126").Value = 5026.1577
$ws.Range("K126").Value = 7462.799999999999
$ws.Range("L126").Value = 15078.4731
$ws.Range("M126").Value = -4992.799999999999
$ws.Range("N126").Value = -20018.4731
$ws.Range("H132").Value = 5000.9414
$ws.Range("I132").Value = 4173.6665
$ws.Range("J132").Value = 5931.625
$ws.Range("K132").Value = 12520.9995
$ws.Range("L132").Value = 17794.875
$ws.Range("M132").Value = -9990.999500000002
$ws.Range("N132").Value = -22854.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4724.3335
$ws.Range("I40").Value = 4283.3687
$ws.Range("K40").Value = 4283.3687
$ws.Range("M40").Value = -4147.3687
$ws.Range("H132").Value = 20011876
$ws.Range("I132").Value = 9748.75
$ws.Range("J132").Value = 29424642
$ws.Range("K132").Value = 29246.25
$ws.Range("L132").Value = 88273926
$ws.Range("M132").Value = -26716.25
$ws.Range("N132").Value = -88278986
$ws.Range("H136").Value = 8988.883
$ws.Range("I136").Value = 2057.5715
$ws.Range("J136").Value = 41335
$ws.Range("K136").Value = 6172.7145
$ws.Range("L136").Value = 124005
$ws.Range("M136").Value = -3622.7145
$ws.Range("N136").Value = -129105

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2569.1304
$ws.Range("I126").Value = 1824.5
$ws.Range("K126").Value = 5473.5
$ws.Range("M126").Value = -3003.5
$ws.Range("H132").Value = 2562
$ws.Range("I132").Value = 1969.1818
$ws.Range("J132").Value = 3866.2
$ws.Range("K132").Value = 5907.5454
$ws.Range("L132").Value = 11598.6
$ws.Range("M132").Value = -3377.5454
$ws.Range("N132").Value = -16658.6
$ws.Range("H136").Value = 1592.8
$ws.Range("I136").Value = 1800.1818
$ws.Range("J136").Value = 1339.3334
$ws.Range("K136").Value = 5400.5454
$ws.Range("L136").Value = 4018.0002
$ws.Range("M136").Value = -2850.5454
$ws.Range("N136").Value = -9118.0002
